# daily auto push: 2025-10-09 22:29 UTC
# Append the next timing-log entry (row 87) to the bottom of the sheet's
# data table: date 2025/10/10, weekday 金 (Friday), time 6, ranking 30.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A stores the date as plain text (matching every other row in the
# sheet, e.g. A86 == "2025/10/10"), not as a native Excel date serial.
# A bare `.Value = "2025/10/10"` would be auto-recognized as a date by
# Excel's input parser, so force the cell to text first, then restore the
# default (General) formatting once the literal string is safely stored.
$ws.Range("A87").NumberFormat = "@"
$ws.Range("A87").Value = "2025/10/10"
$ws.Range("A87").ClearFormats()

$ws.Range("B87").Value = "金"
$ws.Range("C87").Value = 6
$ws.Range("D87").Value = 30
